# Fix bug where the model always output the same prediction:
# refresh rows 2-5 with corrected values and append additional
# prediction rows (6-14 and 18-20; rows 15-17 intentionally absent).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row, delay, gain, perc_trade, number_of_trade
$data = @(
    @(2,  1,  -583.8945067295343,  21.90919926319211,  1.083462132921175,  18458),
    @(3,  2,  -375.2425857536614,  17.02335157392859,  1.17439446366782,   18457),
    @(4,  3,  -309.4941282964456,  12.93346337234504,  1.226679104477612,  18456),
    @(5,  4,  -235.2131625981372,  9.536710918450284,  1.205513784461153,  18455),
    @(6,  5,  -287.6495036159624,  9.846103825728839,  1.057757644394111,  18454),
    @(7,  6,  -247.434345068912,   8.822413699669431,  1.167776298268975,  18453),
    @(8,  7,  -157.8501286614251,  8.351398222414915,  1.269513991163476,  18452),
    @(9,  8,  -219.0863649459424,  8.124220909435802,  0.9646133682830931, 18451),
    @(10, 9,  -116.822503001472,   6.699186991869918,  1.243194192377495,  18450),
    @(11, 10, -179.1598834748678,  6.601983847363001,  1.096385542168675,  18449),
    @(12, 11, -292.8185644338039,  6.526452732003469,  0.8871473354231975, 18448),
    @(13, 12, -166.9050034017448,  5.849189570119803,  1.111545988258317,  18447),
    @(14, 13, -218.4570402655779,  5.665184863927139,  0.9387755102040817, 18446),
    @(18, 17, -228.4948665390954,  4.777139138922026,  0.8430962343096234, 18442),
    @(19, 18, -241.9085060696412,  4.576758310286861,  0.7805907172995781, 18441),
    @(20, 19, 60.8864454667772,    4.202819956616052,  1.327327327327327,  18440)
)

foreach ($row in $data) {
    $r = $row[0]
    $values = New-Object 'object[,]' 1,5
    $values[0,0] = $row[1]
    $values[0,1] = $row[2]
    $values[0,2] = $row[3]
    $values[0,3] = $row[4]
    $values[0,4] = $row[5]
    $ws.Range("A${r}:E${r}").Value = $values
}

Write-Output "Updated prediction rows; dimension now A1:E20"
